$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - prefix with apostrophe to force text
# (values look numeric, e.g. "313.83", and Excel would otherwise parse them
# as numbers), then reset the style back to Normal so no stray "quotePrefix"
# style gets introduced (matches original formatting).
$ws.Range("D2").Value = "'27.386.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.824.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Value = "'313.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'0.4664"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.3791"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.07449"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.8758"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'20.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'1.826.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'6.684"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'5.422"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'93.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.07090"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "'0.000008794"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'27.384.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'5.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.052.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Value = "'151.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Value = "'18.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'5.340"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'117.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'0.08970"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.7889"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Value = "'4.539"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'2.944"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'0.9997"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'1.100"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "'0.01976"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'0.05252"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Value = "'7.306"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.5373"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'2.352"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.1703"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'8.667"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.5104"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'10.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'105.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'1.683"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.9995"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.06383"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values - plain percentage text with padding spaces
$ws.Range("E2").Value = "  +2.75%  "
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  +4.80%  "
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("E12").Value = "  -7.24%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +3.06%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +2.59%  "
$ws.Range("E22").Value = "  +3.57%  "
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  +3.80%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +6.71%  "
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("E34").Value = "  +2.44%  "
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +4.34%  "
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("E43").Value = "  +20.53%  "
$ws.Range("E44").Value = "  +1.74%  "
$ws.Range("E45").Value = "  +3.47%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  +1.46%  "
